# Revert "Merge remote-tracking branch 'origin/Dev_0.0.1' into ArtWork"
# Rewrites the #Spwaner sheet back to its pre-merge content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 1 (header / field names)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "DistanceMin"
$ws.Range("C1").Value = "DistanceMax"
$ws.Range("D1").Value = "RepeatDistance"
$ws.Range("E1").Value = "Rate"
$ws.Range("F1").Value = "MaxSpwan"
$ws.Range("G1").Value = "Monster"
$ws.Range("H1").Value = "Score"
$ws.Range("I1").Value = "IncreaseAtkScale"
$ws.Range("J1").Value = "IncreaseHpScale"
$ws.Range("K1").Value = "Lood"
$ws.Range("L1").Value = "MoveKey"
$ws.Range("M1").ClearContents()

# ---------------------------------------------------------------
# Row 2 (field descriptions)
# ---------------------------------------------------------------
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "해당 거리부터 등장"
$ws.Range("C2").Value = "해당 거리까지 등장"
$ws.Range("D2").Value = "해당 거리마다 스폰실행"
$ws.Range("E2").Value = "스폰 확률"
$ws.Range("F2").Value = "최대 스폰 횟수`n0 = 제한없음"
$ws.Range("G2").Value = "등장시킬 몬스터"
$ws.Range("H2").Value = "점수"
$ws.Range("I2").Value = "공격력 증가 배율"
$ws.Range("J2").Value = "체력 증가 배율"
$ws.Range("K2").Value = "등장 좌표"
$ws.Range("L2").Value = "이동패턴"
$ws.Range("M2").ClearContents()

# ---------------------------------------------------------------
# Row 3 (field types)
# ---------------------------------------------------------------
$ws.Range("A3").Value = "string"
$ws.Range("B3").Value = "long"
$ws.Range("C3").Value = "long"
$ws.Range("D3").Value = "long"
$ws.Range("E3").Value = "long"
$ws.Range("F3").Value = "long"
$ws.Range("G3").Value = "string"
$ws.Range("H3").Value = "long"
$ws.Range("I3").Value = "float"
$ws.Range("J3").Value = "float"
$ws.Range("K3").Value = "string"
$ws.Range("L3").Value = "string"
$ws.Range("M3").ClearContents()

# ---------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------
$ws.Range("A4").Value = "잼민이"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1000
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "잼민이"
$ws.Range("H4").Value = 10
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = "0,0,0"
$ws.Range("L4").Value = "left"
$ws.Range("M4").ClearContents()

# ---------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------
$ws.Range("A5").Value = "급식충"
$ws.Range("B5").Value = 1000
$ws.Range("C5").Value = 3000
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = 0.8
$ws.Range("F5").Value = 0
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").Value = "right"
$ws.Range("M5").ClearContents()

# ---------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------
$ws.Range("A6").Value = "학식충"
$ws.Range("B6").Value = 2000
$ws.Range("C6").Value = 4000
$ws.Range("D6").Value = 200
$ws.Range("E6").Value = 0.9
$ws.Range("F6").Value = 0
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()

# ---------------------------------------------------------------
# Row 7
# ---------------------------------------------------------------
$ws.Range("A7").Value = "금태양"
$ws.Range("B7").Value = 2000
$ws.Range("C7").Value = 4000
$ws.Range("D7").Value = 500
$ws.Range("E7").Value = 0.02
$ws.Range("F7").Value = 2
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()

# ---------------------------------------------------------------
# Row 8 (new row)
# ---------------------------------------------------------------
$ws.Range("B8").Value = 10000
$ws.Range("C8").Value = "max"
$ws.Range("D8").Value = 200

# ---------------------------------------------------------------
# Drop the now-unused column M entirely
# ---------------------------------------------------------------
$ws.Columns.Item(13).Delete()

# ---------------------------------------------------------------
# Restore the previously-selected cell
# ---------------------------------------------------------------
$ws.Range("K4").Select()
